# Apply the latest scraped crypto price/volume snapshot to Sheet1.
# Cells whose text looks like a plain number (e.g. "22.74") are written
# with a leading apostrophe so Excel stores them as TEXT (quote-prefixed),
# matching the original inline-string cells instead of being reinterpreted
# as numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.132.76'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '1.677.90'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range("D5").Value = "`'214.13"
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range("D8").Value = "`'22.74"
$ws.Range('E8').Value = '  +6.73%  '
$ws.Range('E9').Value = '  +2.17%  '
$ws.Range("D10").Value = "`'0.0621"
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').Value = '1.915.77'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '1.680.53'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range("D14").Value = "`'4.20"
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range("D15").Value = "`'0.550"
$ws.Range('E15').Value = '  +3.16%  '
$ws.Range("D16").Value = "`'66.52"
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '27.101.60'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range("D18").Value = "`'235.41"
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range("D19").Value = "`'7.88"
$ws.Range('E19').Value = '  -3.15%  '
$ws.Range('D20').Value = '0.0₃0740'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range("D22").Value = "`'4.53"
$ws.Range('E22').Value = '  +1.60%  '
$ws.Range("D23").Value = "`'9.50"
$ws.Range('E23').Value = '  +2.52%  '
$ws.Range('E24').Value = '  -1.04%  '
$ws.Range("D25").Value = "`'147.42"
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range("D26").Value = "`'7.41"
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range("D27").Value = "`'16.31"
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +0.66%  '
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').Value = '1.541.56'
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range("D34").Value = "`'3.23"
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('E35').Value = '  -3.11%  '
$ws.Range("D36").Value = "`'0.606"
$ws.Range('E36').Value = '  +3.11%  '
$ws.Range('E37').Value = '  +2.39%  '
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('E40').Value = '  +2.72%  '
$ws.Range('E41').Value = '  +3.87%  '
$ws.Range("D42").Value = "`'69.41"
$ws.Range('E42').Value = '  +2.09%  '
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('D45').Value = '1.822.62'
$ws.Range("D46").Value = "`'0.778"
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range("D47").Value = "`'89.95"
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = "`'1.63"
$ws.Range('E48').Value = '  +6.18%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0111'
$ws.Range('E49').Value = '  +3.07%  '
$ws.Range("D50").Value = "`'8.20"
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('E51').Value = '  -0.20%  '
